$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R1")

# Row 4 loses its outage details but keeps the hub-site identifier (D4).
$ws.Range("B4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("L4").Value = ""

# Rows 5-7 (the other R4 outage rows) are removed entirely, shrinking the
# used range down to A1:L4.
$ws.Range("A5:L7").Delete()
